# Add a new trailing period column (BB) to the quarterly series table.
# Column BB mirrors column BA (the latest period's YoY figures get
# "carried forward") and BB1 gets the new period-end date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new date in BB1, using the same style/number format as BA1 ---
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BB1").Value = 45986
$excel.CutCopyMode = 0

# --- Data rows: carry the last available value (column BA) into column BB ---
for ($r = 3; $r -le 21; $r++) {
  $srcCell = $ws.Cells.Item($r, 53)   # column BA
  $dstCell = $ws.Cells.Item($r, 54)   # column BB
  $srcVal = $srcCell.Value()
  if ($null -ne $srcVal) {
    $dstCell.Value = $srcVal
  }
}
